$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) First paragraph: add two trailing spaces to the existing text,
#    then append three red-colored runs forming:
#    "(This is a change – Version for main branch)"
# -----------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "This is a Microsoft word document.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This is a Microsoft word document.  ", 2)

$p1 = $d.Paragraphs.Item(1)
$p1End = $p1.Range.End
$insertPos = $p1End - 1   # position just before the paragraph mark

function Add-RedRun([int]$pos, [string]$text) {
    $ip = $word.ActiveDocument.Range($pos, $pos)
    $ip.InsertAfter($text)
    $newEnd = $pos + $text.Length
    $coloredRange = $word.ActiveDocument.Range($pos, $newEnd)
    $coloredRange.Font.Color = 255
    return $newEnd
}

$insertPos = Add-RedRun $insertPos "(This is a change – Ve"
$insertPos = Add-RedRun $insertPos "rsion for main branch"
$insertPos = Add-RedRun $insertPos ")"

# -----------------------------------------------------------------
# 2) Append a new, empty paragraph after the very last paragraph in
#    the document body, shaded with fill color F9F9F9.
#    InsertXML is used so the paragraph is created clean (no
#    inherited run/paragraph formatting from the preceding
#    paragraph).
# -----------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>'
$endRange.InsertXML($newParaXml)

# -----------------------------------------------------------------
# 3) Remove now-unused built-in/custom styles (Heading 2/4 and their
#    linked char styles, Hyperlink, and various custom styles that
#    are not referenced anywhere in the document body). Deletion must
#    proceed from the highest style index to the lowest to avoid
#    invalidating earlier indices.
# -----------------------------------------------------------------
$styleIndicesToDelete = @(18, 17, 16, 15, 14, 13, 12, 11, 10, 3, 2)
foreach ($idx in $styleIndicesToDelete) {
    $style = $d.Styles.Item($idx)
    $style.Delete()
}

Write-Host "Edit complete"
